$wb = $excel.ActiveWorkbook

# --- Sheet "Tactile Tabletop Data - Level 2": add two new card rows (42 & 43) ---
$ws1 = $wb.Worksheets.Item("Tactile Tabletop Data - Level 2")

# Row 42 - "Exploitation" card
$ws1.Range("A42").Value = "Exploitation"
$ws1.Range("B42").Value = "Self"
$ws1.Range("C42").Value = "Instant"
$ws1.Range("F42").Value = "You may play one addition cards top or bottom action in place of this action"
$ws1.Range("G42").Value = "Discard"
$ws1.Range("H42").Value = "Abuse Knowledge"
$ws1.Range("I42").Value = "Self"
$ws1.Range("J42").Value = "Instant"
$ws1.Range("M42").Value = "You may choose the top and bottom action of one card you played this turn"
$ws1.Range("N42").Value = "Exhaust"
$ws1.Range("O42").Value = "2 Level Points"
$ws1.Range("P42").Value = "2 Kno"
$ws1.Rows.Item(42).RowHeight = 38.25

# Row 43 - "Shelter Friends" card
$ws1.Range("A43").Value = "Shelter Friends"
$ws1.Range("B43").Value = "Allies"
$ws1.Range("C43").Value = "1 RND"
$ws1.Range("D43").Value = "X = Influence"
$ws1.Range("F43").Value = "For all adjacent allies, if they take damage, instead you take that damage and reduce it by X"
$ws1.Range("G43").Value = "Exhaust"
$ws1.Range("H43").Value = "Healing Sap"
$ws1.Range("I43").Value = "Ally"
$ws1.Range("J43").Value = "Instant"
$ws1.Range("K43").Value = "X = Influence    Y = Level"
$ws1.Range("M43").Value = "Heal Target Ally for Y and lose X life"
$ws1.Range("N43").Value = "Discard"
$ws1.Range("O43").Value = "2 Level Points"
$ws1.Range("P43").Value = "2 Vig, 1 Spi"
$ws1.Rows.Item(43).RowHeight = 38.25

# --- Sheet "Sheet1": mark three existing draft rows as DONE (new column I) ---
$ws2 = $wb.Worksheets.Item("Sheet1")

$ws2.Range("I6").Value = "DONE"
$ws2.Range("I6").WrapText = $true

$ws2.Range("I10").Value = "DONE"
$ws2.Range("I10").WrapText = $true

$ws2.Range("I11").Value = "DONE"
$ws2.Range("I11").WrapText = $true

# --- View state: reflect the user's final selection (Sheet1 selection first,
#     then finish on the data sheet so it ends up the active/selected tab) ---
$ws2.Activate() | Out-Null
$ws2.Range("S31:T34").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("F50").Select() | Out-Null
